# Update LR-pairs TPM values for Rspo3-Fzd8 sheet with newly computed stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 0.9300694554254023
$ws.Range("J2").Value = 0.9300694554254023
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.083576666666666
$ws.Range("N2").Value = 9.250729999999999
$ws.Range("O2").Value = 0.2272509363535097
$ws.Range("P2").Value = 0.2272509363535097
$ws.Range("Q2").Value = 7.732442632302221
$ws.Range("R2").Value = 69.59198369072
$ws.Range("S2").Value = 0.2113591546192215
$ws.Range("T2").Value = 0.2113591546192215

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.9300694554254023
$ws.Range("J3").Value = 0.9300694554254023
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("O3").Value = 0.4756405360586227
$ws.Range("P3").Value = 0.4756405360586227
$ws.Range("S3").Value = 0.4423787343502897
$ws.Range("T3").Value = 0.4423787343502896

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.9300694554254023
$ws.Range("J4").Value = 0.9300694554254023
$ws.Range("M4").Value = 4.031477000000001
$ws.Range("N4").Value = 12.094431
$ws.Range("O4").Value = 0.2971085275878677
$ws.Range("P4").Value = 0.2971085275878677
$ws.Range("Q4").Value = 10.10941773004267
$ws.Range("R4").Value = 90.98475957038401
$ws.Range("S4").Value = 0.2763315664558912
$ws.Range("T4").Value = 0.2763315664558912

# Row 5 (MuSCs -> ECs)
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1885443333333333
$ws.Range("H5").Value = 0.5656329999999999
$ws.Range("I5").Value = 0.06993054457459773
$ws.Range("J5").Value = 0.06993054457459771
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.083576666666666
$ws.Range("N5").Value = 9.250729999999999
$ws.Range("O5").Value = 0.2272509363535097
$ws.Range("P5").Value = 0.2272509363535097
$ws.Range("Q5").Value = 0.5813909068988887
$ws.Range("R5").Value = 5.232518162089999
$ws.Range("S5").Value = 0.01589178173428818
$ws.Range("T5").Value = 0.01589178173428818

# Row 6 (MuSCs -> FAPs)
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1885443333333333
$ws.Range("H6").Value = 0.5656329999999999
$ws.Range("I6").Value = 0.06993054457459773
$ws.Range("J6").Value = 0.06993054457459771
$ws.Range("M6").Value = 6.453984666666667
$ws.Range("O6").Value = 0.4756405360586227
$ws.Range("P6").Value = 0.4756405360586227
$ws.Range("Q6").Value = 1.216862236320222
$ws.Range("R6").Value = 10.951760126882
$ws.Range("S6").Value = 0.03326180170833307
$ws.Range("T6").Value = 0.03326180170833307

# Row 7 (MuSCs -> MuSCs)
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1885443333333333
$ws.Range("H7").Value = 0.5656329999999999
$ws.Range("I7").Value = 0.06993054457459773
$ws.Range("J7").Value = 0.06993054457459771
$ws.Range("M7").Value = 4.031477000000001
$ws.Range("N7").Value = 12.094431
$ws.Range("O7").Value = 0.2971085275878677
$ws.Range("P7").Value = 0.2971085275878677
$ws.Range("Q7").Value = 0.7601121433136667
$ws.Range("R7").Value = 6.841009289823
$ws.Range("S7").Value = 0.02077696113197648
$ws.Range("T7").Value = 0.02077696113197647
